$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume-change (E) figures, and the three-row
# reorderings in B/C/D/E for rows 40-41 and 49-51, per the latest crypto pull.

$ws.Range("D2").Value = "99.178.55"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.293.79"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'255.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'625.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +23.70%  "
$ws.Range("D8").Value = "'0.403"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.21%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.982"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +24.33%  "
$ws.Range("D11").Value = "3.291.67"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "'40.73"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +13.74%  "
$ws.Range("D14").Value = "98.843.11"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "'0.0000250"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "3.912.07"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "'5.46"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "3.291.31"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").Value = "'6.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +10.45%  "
$ws.Range("D21").Value = "'15.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").Value = "'492.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.37%  "
$ws.Range("D23").Value = "'9.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").Value = "'0.0000203"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'0.345"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +41.56%  "
$ws.Range("D26").Value = "'5.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").Value = "'89.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "'12.16"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").Value = "3.467.94"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'0.140"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +14.94%  "
$ws.Range("D32").Value = "'0.190"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'10.70"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +16.57%  "
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "'28.03"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  +7.47%  "
$ws.Range("D37").Value = "'0.151"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'7.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").Value = "'1.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'492.57"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'24.71"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "'3.82"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.94%  "
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "'3.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("D47").Value = "'158.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'7.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +16.84%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.78"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.97%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.851"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.85%  "
